# Updated cryptos list on Mon Oct  2 10:43:10 UTC 2023 with GitHub Actions
#
# Price (D) and Volume(1h) (E) cells are stored as literal text in this sheet (e.g. "28.309.53",
# "  +4.00%  "), not real numbers. For Price values that look like a plain decimal number
# (e.g. "219.35", "70.30"), a bare Range.Value assignment would let Excel's normal General-format
# autoconvert turn them into real numbers -- silently retyping the cell and, for values with an
# insignificant trailing zero such as "70.30" / "2.30", collapsing it to 70.3 / 2.3. Prefixing the
# literal with an apostrophe (the PowerShell '' is an escaped single quote, i.e. the text fed to Excel
# starts with ') forces Excel to keep the entry as text, exactly like a user typing `'70.30` into the
# cell, which preserves both the exact digits and the original Text cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.309.53'
$ws.Range('E2').Value = '  +4.00%  '

$ws.Range('D3').Value = '''1.731.66'
$ws.Range('E3').Value = '  +2.68%  '

$ws.Range('D5').Value = '''219.35'
$ws.Range('E5').Value = '  +1.57%  '

$ws.Range('D6').Value = '''0.523'
$ws.Range('E6').Value = '  +0.56%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').Value = '''24.19'
$ws.Range('E8').Value = '  +6.32%  '

$ws.Range('E9').Value = '  +2.99%  '

$ws.Range('E10').Value = '  +1.52%  '

$ws.Range('D11').Value = '''0.0894'
$ws.Range('E11').Value = '  +0.39%  '

$ws.Range('D12').Value = '''1.977.22'
$ws.Range('E12').Value = '  +2.73%  '

$ws.Range('D13').Value = '''1.732.34'
$ws.Range('E13').Value = '  +2.81%  '

$ws.Range('D14').Value = '''4.26'
$ws.Range('E14').Value = '  +1.94%  '

$ws.Range('E15').Value = '  +1.82%  '

$ws.Range('D16').Value = '''67.78'
$ws.Range('E16').Value = '  +0.86%  '

$ws.Range('D17').Value = '''28.303.63'
$ws.Range('E17').Value = '  +3.99%  '

$ws.Range('D18').Value = '''243.37'
$ws.Range('E18').Value = '  +1.81%  '

$ws.Range('E19').Value = '  +1.34%  '

$ws.Range('E20').Value = '  -2.66%  '

$ws.Range('E21').Value = '  -0.14%  '

$ws.Range('D22').Value = '''4.66'
$ws.Range('E22').Value = '  +2.03%  '

$ws.Range('D23').Value = '''9.77'
$ws.Range('E23').Value = '  +1.51%  '

$ws.Range('E24').Value = '  -0.41%  '

$ws.Range('D25').Value = '''149.38'
$ws.Range('E25').Value = '  +0.73%  '

$ws.Range('E26').Value = '  +3.24%  '

$ws.Range('D27').Value = '''16.63'
$ws.Range('E27').Value = '  +0.82%  '

$ws.Range('E28').Value = '  +0.95%  '

$ws.Range('E29').Value = '  -0.04%  '

$ws.Range('D30').Value = '''0.0517'
$ws.Range('E30').Value = '  +3.11%  '

$ws.Range('E31').Value = '  +2.25%  '

$ws.Range('E32').Value = '  +0.56%  '

$ws.Range('D33').Value = '''3.28'

$ws.Range('D34').Value = '''1.488.77'
$ws.Range('E34').Value = '  -5.50%  '

$ws.Range('E35').Value = '  -1.94%  '

$ws.Range('D36').Value = '''0.979'
$ws.Range('E36').Value = '  +2.52%  '

$ws.Range('D37').Value = '''0.603'
$ws.Range('E37').Value = '  +0.08%  '

$ws.Range('E38').Value = '  +0.53%  '

$ws.Range('E39').Value = '  +1.08%  '

$ws.Range('E40').Value = '  +0.51%  '

$ws.Range('D41').Value = '''70.30'
$ws.Range('E41').Value = '  +0.73%  '

$ws.Range('E42').Value = '  -0.07%  '

$ws.Range('E43').Value = '  +0.30%  '

$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '''2.30'
$ws.Range('E44').Value = '  +1.88%  '

$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '''1.880.92'
$ws.Range('E45').Value = '  +2.48%  '

$ws.Range('D46').Value = '''0.796'
$ws.Range('E46').Value = '  +1.22%  '

$ws.Range('D47').Value = '''1.73'
$ws.Range('E47').Value = '  +7.78%  '

$ws.Range('E48').Value = '  +5.99%  '

$ws.Range('D49').Value = '''90.84'
$ws.Range('E49').Value = '  -0.50%  '

$ws.Range('D50').Value = '''8.25'
$ws.Range('E50').Value = '  +0.23%  '

$ws.Range('D51').Value = '''0.104'
$ws.Range('E51').Value = '  -0.48%  '
